$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestReportTest")

# Row 4: change Run Status from "run" to "skip", and update the Params text
$ws.Range("B4").Value = "skip"
$ws.Range("C4").Value = "Automation_Test,Vibhor,model,Updated Description,Test Description,vivek"

# New row 5: negativescenerios_testreport / run / Automation_Test,Vibhor,model,Test.jpeg,Test Description
$ws.Range("A5").Value = "negativescenerios_testreport"
$ws.Range("B5").Value = "run"
$ws.Range("C5").Value = "Automation_Test,Vibhor,model,Test.jpeg,Test Description"

# Widen column C to fit the longer text
$ws.Columns.Item(3).ColumnWidth = 59.2421875

# Update the active selection to match the authored state
$ws.Range("C16").Select()
